$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): B1/C1 content swap (FirstName/LastName) ---
# A1 keeps "CNE" but gets re-touched with the (effectively default) style.
$ws.Range("A1").Value = "CNE"
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"

# Re-apply the (default/"Normal") style on the header's first three cells so
# they end up referencing the newly introduced cellXf (identical to the
# workbook's default formatting) instead of the bold/explicit-color one.
$ws.Range("A1:C1").Style = "Normal"

# --- Data rows 2-11: update the CNE id in column A and touch style on A:C ---
$ids = @(18000001, 18000002, 18000003, 18000004, 18000005, 18000006, 18000007, 18000008, 18000009, 18000010)
for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $ids[$i]
    $ws.Range("A" + $row + ":C" + $row).Style = "Normal"
}

# --- Restore the active selection to match the saved workbook state ---
$ws.Range("F14").Select()
